{"js": "// The \"Bibliografia\" paragraph holds a single run whose text is six\n// numbered references glued together with no separators. Split it into\n// six <w:t> runs of text joined by manual line breaks (<w:br/>), all\n// still inside the same run, matching the target OOXML.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker =\n  \"1)Telles, P. C. S. - Materiais para Equipamentos de Processo\";\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text && p.text.indexOf(marker) !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate the Bibliografia paragraph.\");\n}\n\nconst parts = [\n  \"1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interci\u00eancia, 4\u00ba Ed., 1989.\",\n  \"2)Bresciani, F., E. - Sele\u00e7\u00e3o de Materiais Met\u00e1licos - Ed. da UNICAMP, 2\u00ba Ed.\",\n  \"3)Freire, J. M. -Materiais de Constru\u00e7\u00e3o Mec\u00e2nica - Ed. Livros T\u00e9cnicos e Cient\u00edficos, Editora 1993.\",\n  \"4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2\u00aa Edi\u00e7\u00e3o.\",\n  \"5)Chiaverini, V.Tecnologia Mec\u00e2nica - Materiais de Constru\u00e7\u00e3o Mec\u00e2nica - Vol. II - Ed. McGraw Hill do Brasil Ltda.\",\n  \"6)Gentil, V. - Corros\u00e3o. - Ed. Guanabara Dois, 1982.\",\n];\n\n// U+000B (vertical tab) is how Word's text model represents a manual\n// line break (<w:br/>) inline within a run's text content.\nconst joined = parts.join(\"\\x0b\");\n\nconst range = target.getRange(\"Whole\");\nrange.insertText(joined, \"Replace\");\n\nawait context.sync();\n", "ps1": "# The \"Bibliografia\" paragraph holds a single run whose text is six\n# numbered references glued together with no separators. Split it into\n# six lines joined by manual line breaks (Chr(11), i.e. vertical tab),\n# which Word renders/serializes as <w:br/> inside the same run.\n\n$d = $word.ActiveDocument\n\n$parts = @(\n  \"1)Telles, P. C. S. - Materiais para Equipamentos de Processo - Ed. Interci\u00eancia, 4\u00ba Ed., 1989.\",\n  \"2)Bresciani, F., E. - Sele\u00e7\u00e3o de Materiais Met\u00e1licos - Ed. da UNICAMP, 2\u00ba Ed.\",\n  \"3)Freire, J. M. -Materiais de Constru\u00e7\u00e3o Mec\u00e2nica - Ed. Livros T\u00e9cnicos e Cient\u00edficos, Editora 1993.\",\n  \"4)A. Remy/ M. Gay/ R. Gonthier - Materiais - Hemus Editora Limitada - 2\u00aa Edi\u00e7\u00e3o.\",\n  \"5)Chiaverini, V.Tecnologia Mec\u00e2nica - Materiais de Constru\u00e7\u00e3o Mec\u00e2nica - Vol. II - Ed. McGraw Hill do Brasil Ltda.\",\n  \"6)Gentil, V. - Corros\u00e3o. - Ed. Guanabara Dois, 1982.\"\n)\n\n$joined = [string]::Join([string][char]11, $parts)\n\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*1)Telles, P. C. S.*\") {\n        $p.Range.Text = $joined\n        break\n    }\n}\n"}
